$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("2025")
$ws.Range("E2").Value = 29945.71989099024
$ws.Range("G2").Value = 8095.925712661654
$ws.Range("I2").Value = 13698.12725754988
$ws.Range("L2").Value = 54033.78794259601
$ws.Range("M2").Value = 10556.13095757
$ws.Range("N2").Value = 7648.163036096154
$ws.Range("O2").Value = 7602.027316031194

$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 4291.832072666659
$ws.Range("E2").Value = 58123.2420465009
$ws.Range("G2").Value = 8095.925712661654
$ws.Range("I2").Value = 27170.44055490107
$ws.Range("L2").Value = 95093.05349626098
$ws.Range("M2").Value = 20633.48636216725
$ws.Range("N2").Value = 10954.02563372962
$ws.Range("O2").Value = 9635.999202966375

$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 1524.327458338266
$ws.Range("B2").Value = 7110.124845112752
$ws.Range("E2").Value = 71988.56920036966
$ws.Range("G2").Value = 8095.925712661654
$ws.Range("I2").Value = 44419.89816206333
$ws.Range("L2").Value = 95093.05349626098
$ws.Range("M2").Value = 23672.65274446728
$ws.Range("N2").Value = 15993.02254519105
$ws.Range("O2").Value = 15252.03544981399

$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 1524.327458338266
$ws.Range("B2").Value = 7110.124845112752
$ws.Range("E2").Value = 71988.56920036966
$ws.Range("G2").Value = 8095.925712661654
$ws.Range("I2").Value = 44419.89816206333
$ws.Range("L2").Value = 95093.05349626098
$ws.Range("M2").Value = 23672.65274446728
$ws.Range("N2").Value = 15993.02254519105
$ws.Range("O2").Value = 15252.03544981399

$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 1524.327458338266
$ws.Range("B2").Value = 7110.124845112752
$ws.Range("E2").Value = 71988.56920036966
$ws.Range("G2").Value = 8095.925712661654
$ws.Range("I2").Value = 44419.89816206333
$ws.Range("L2").Value = 95093.05349626098
$ws.Range("M2").Value = 23672.65274446728
$ws.Range("N2").Value = 15993.02254519105
$ws.Range("O2").Value = 15252.03544981399

$ws = $wb.Worksheets.Item("2050")
$ws.Range("A2").Value = 1524.327458338266
$ws.Range("B2").Value = 7110.124845112752
$ws.Range("E2").Value = 71988.56920036966
$ws.Range("G2").Value = 8095.925712661654
$ws.Range("I2").Value = 44419.89816206333
$ws.Range("L2").Value = 95093.05349626098
$ws.Range("M2").Value = 23672.65274446728
$ws.Range("N2").Value = 15993.02254519105
$ws.Range("O2").Value = 15252.03544981399

